# Auto - Update data with bot!
# Updates specific blog entries (title/link) in the "Blogs_used_list" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: "연속 신호의 샘플링" -> "C에서 MATLAB 호출하기"
$ws.Range("D5").Value = "C에서 MATLAB 호출하기"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2022/06/20/calling_MATLAB_in_C.html"

# Row 26: "ai plus(est soft)" -> "인공지능 음성 생성 연구: 음성 분류 솔루션"
$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

# Row 37: Paper review title update + link uid change
$ws.Range("D37").Value = "[Paper Review] Unsupervised Keyphrase Extraction by Jointly Modeling Local and Global Context"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=2029&mod=document&pageid=1"

# Row 41: "DevOps 에 대한 이해" -> "cloudinsight"
$ws.Range("D41").Value = "cloudinsight"

# Row 51: vim title -> pandas title + link number change
$ws.Range("D51").Value = "[python + pandas] 데이터프레임에서 특정 기간의 데이터 추출하기"
$ws.Range("E51").Value = "https://bskyvision.com/1304"
